# Refresh the cryptocurrency price (column D) and 1h volume-change
# (column E) figures for the rows whose upstream data moved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores figures as literal text (e.g. "1.827.81",
# "0.9999") rather than numbers, so force the Text format before writing
# the new values -- otherwise Excel's automatic "looks like a number"
# detection would silently convert values such as "311.29" into the
# numeric 311.29 and drop formatting like the trailing zero in "1.000".
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.042.49"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "1.822.46"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -1.18%  "
$ws.Range("D5").Value = "311.29"
$ws.Range("E5").Value = "  -2.34%  "
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("D7").Value = "0.4223"
$ws.Range("E7").Value = "  -1.64%  "
$ws.Range("D8").Value = "0.3669"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "0.07214"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").Value = "0.8422"
$ws.Range("E10").Value = "  -3.80%  "
$ws.Range("D11").Value = "20.75"
$ws.Range("E11").Value = "  -3.68%  "
$ws.Range("D12").Value = "1.833.04"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "6.674"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "0.07066"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "5.292"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("D16").Value = "90.31"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "0.000008766"
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").Value = "14.93"
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("D21").Value = "27.138.08"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").Value = "5.140"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").Value = "2.049.12"
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("D25").Value = "1.980"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "152.01"
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("E27").Value = "  +3.39%  "
$ws.Range("D28").Value = "18.25"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").Value = "5.271"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").Value = "117.57"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").Value = "0.08731"
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("D32").Value = "1.178"
$ws.Range("E32").Value = "  -3.81%  "
$ws.Range("D33").Value = "0.7378"
$ws.Range("E33").Value = "  -4.94%  "
$ws.Range("D34").Value = "2.903"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").Value = "4.420"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("D37").Value = "1.091"
$ws.Range("E37").Value = "  -3.36%  "
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").Value = "0.05258"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").Value = "7.330"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").Value = "2.873"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("D42").Value = "0.1687"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "0.5045"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").Value = "8.564"
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("D45").Value = "10.47"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").Value = "106.20"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").Value = "0.4710"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").Value = "1.917"
$ws.Range("E48").Value = "  +4.49%  "
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "0.06336"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("D51").Value = "1.649"
$ws.Range("E51").Value = "  -2.15%  "
